# Add data for 2023-09-04
# Update the 2023 (column J) violent-crime totals on the citywide summary
# sheet, the "By Neighborhood" roll-up sheet, and every individual
# neighborhood sheet that has new counts for this refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("J2").Value = 5128
$ws.Range("J3").Value = 5472
$ws.Range("J4").Value = 1208
$ws.Range("J6").Value = 6808
$ws.Range("J7").Value = 19042

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("J6").Value = 163
$ws.Range("J7").Value = 254

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("J3").Value = 302
$ws.Range("J4").Value = 56
$ws.Range("J6").Value = 230
$ws.Range("J7").Value = 816

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("J2").Value = 118
$ws.Range("J3").Value = 150
$ws.Range("J6").Value = 109
$ws.Range("J7").Value = 398

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("J3").Value = 97
$ws.Range("J7").Value = 283

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("J2").Value = 179
$ws.Range("J3").Value = 286
$ws.Range("J7").Value = 736

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("J2").Value = 65
$ws.Range("J6").Value = 47
$ws.Range("J7").Value = 173

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("J4").Value = 69
$ws.Range("J7").Value = 555
$ws.Range("J8").Value = 1210
$ws.Range("J11").Value = 296
$ws.Range("J15").Value = 208
$ws.Range("J16").Value = 74
$ws.Range("J19").Value = 543
$ws.Range("J20").Value = 397
$ws.Range("J21").Value = 52
$ws.Range("J23").Value = 183
$ws.Range("J27").Value = 113
$ws.Range("J29").Value = 1070
$ws.Range("J31").Value = 173
$ws.Range("J33").Value = 867
$ws.Range("J34").Value = 88
$ws.Range("J36").Value = 262
$ws.Range("J37").Value = 599
$ws.Range("J40").Value = 42
$ws.Range("J42").Value = 781
$ws.Range("J43").Value = 159
$ws.Range("J46").Value = 66
$ws.Range("J47").Value = 144
$ws.Range("J48").Value = 218
$ws.Range("J49").Value = 127
$ws.Range("J52").Value = 483
$ws.Range("J53").Value = 254
$ws.Range("J54").Value = 362
$ws.Range("J55").Value = 245
$ws.Range("J57").Value = 81
$ws.Range("J63").Value = 67
$ws.Range("J65").Value = 498
$ws.Range("J67").Value = 736
$ws.Range("J72").Value = 75
$ws.Range("J73").Value = 176
$ws.Range("J76").Value = 271
$ws.Range("J78").Value = 239
$ws.Range("J79").Value = 547
$ws.Range("J80").Value = 30
$ws.Range("J83").Value = 398
$ws.Range("J85").Value = 816
$ws.Range("J89").Value = 239
$ws.Range("J90").Value = 212
$ws.Range("J93").Value = 82
$ws.Range("J95").Value = 283
$ws.Range("J96").Value = 229
$ws.Range("J100").Value = 39
$ws.Range("J101").Value = 19042

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("J6").Value = 173
$ws.Range("J7").Value = 599

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("J2").Value = 217
$ws.Range("J3").Value = 283
$ws.Range("J6").Value = 292
$ws.Range("J7").Value = 867

$ws = $wb.Worksheets.Item('New City')
$ws.Range("J2").Value = 141
$ws.Range("J6").Value = 178
$ws.Range("J7").Value = 498

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("J6").Value = 74
$ws.Range("J7").Value = 127

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("J6").Value = 170
$ws.Range("J7").Value = 362

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("J2").Value = 319
$ws.Range("J3").Value = 366
$ws.Range("J4").Value = 60
$ws.Range("J6").Value = 283
$ws.Range("J7").Value = 1070

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("J3").Value = 159
$ws.Range("J6").Value = 201
$ws.Range("J7").Value = 543

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("J6").Value = 112
$ws.Range("J7").Value = 218

$ws = $wb.Worksheets.Item('River North')
$ws.Range("J3").Value = 56
$ws.Range("J6").Value = 147
$ws.Range("J7").Value = 271

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("J3").Value = 167
$ws.Range("J6").Value = 177
$ws.Range("J7").Value = 555

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("J2").Value = 176
$ws.Range("J3").Value = 156
$ws.Range("J6").Value = 397
$ws.Range("J7").Value = 781

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("J6").Value = 67
$ws.Range("J7").Value = 239

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("J2").Value = 57
$ws.Range("J7").Value = 245

$ws = $wb.Worksheets.Item('Jefferson Park')
$ws.Range("J3").Value = 14
$ws.Range("J7").Value = 66

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("J2").Value = 50
$ws.Range("J7").Value = 183

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("J2").Value = 76
$ws.Range("J7").Value = 239

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range("J6").Value = 34
$ws.Range("J7").Value = 52

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("J6").Value = 151
$ws.Range("J7").Value = 547

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("J3").Value = 137
$ws.Range("J6").Value = 105
$ws.Range("J7").Value = 397

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("J3").Value = 82
$ws.Range("J7").Value = 262

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range("J4").Value = 8
$ws.Range("J7").Value = 82

$ws = $wb.Worksheets.Item('Wrigleyville')
$ws.Range("J6").Value = 21
$ws.Range("J7").Value = 39

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("J2").Value = 94
$ws.Range("J3").Value = 62
$ws.Range("J6").Value = 116
$ws.Range("J7").Value = 296

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range("J6").Value = 31
$ws.Range("J7").Value = 88

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("J3").Value = 38
$ws.Range("J7").Value = 144

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("J3").Value = 51
$ws.Range("J7").Value = 208

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("J3").Value = 144
$ws.Range("J7").Value = 483

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("J2").Value = 62
$ws.Range("J3").Value = 49
$ws.Range("J7").Value = 176

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("J2").Value = 33
$ws.Range("J3").Value = 25
$ws.Range("J7").Value = 113

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("J3").Value = 61
$ws.Range("J7").Value = 212

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range("J3").Value = 22
$ws.Range("J7").Value = 81

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("J3").Value = 31
$ws.Range("J6").Value = 92
$ws.Range("J7").Value = 159

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("J2").Value = 68
$ws.Range("J7").Value = 229

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("J3").Value = 23
$ws.Range("J7").Value = 75

$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Range("J4").Value = 6
$ws.Range("J7").Value = 30

$ws = $wb.Worksheets.Item('Hegewisch')
$ws.Range("J3").Value = 15
$ws.Range("J7").Value = 42

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("J2").Value = 338
$ws.Range("J3").Value = 367
$ws.Range("J7").Value = 1210

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Range("J2").Value = 23
$ws.Range("J7").Value = 69

$ws = $wb.Worksheets.Item('Bucktown')
$ws.Range("J2").Value = 10
$ws.Range("J7").Value = 74
